$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "Comment" header column (F)
$ws.Range("F1").Value = "Comment"
$ws.Range("F1").Font.Bold = $true

# Row 3 (Guinea Pig) gets a comment equal to the literal text "Comment"
$ws.Range("F3").Value = "Comment"

# Row 4 (Tony Stark) gets a real comment
$ws.Range("F4").Value = "Available comment"
